$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 498, shifting existing rows 498:567 down to 499:568
$ws.Rows(498).Insert()

# Populate the newly inserted row 498 with the new data record
$ws.Range("A498").Value = 5
$ws.Range("B498").Value = "Macroferia Regional de Talca"
$ws.Range("C498").Value = "Maule"
$ws.Range("D498").Value = 45131
$ws.Range("E498").Value = 7
$ws.Range("F498").Value = 100114013
$ws.Range("G498").Value = "Zanahoria"
$ws.Range("H498").Value = "Sin especificar"
$ws.Range("I498").Value = "Primera"
$ws.Range("J498").Value = 600
$ws.Range("K498").Value = 5000
$ws.Range("L498").Value = 5000
$ws.Range("M498").Value = 5000
$ws.Range("N498").Value = "`$/saco 20 kilos"
$ws.Range("O498").Value = "Región de Ñuble"
$ws.Range("P498").Value = 250
$ws.Range("Q498").Value = 20
$ws.Range("R498").Value = "Hortaliza"
